$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Comercializadora del Agro de
# Limarí" Choclo (Dulce o Americano). Insert it at row 11 - this shifts
# all the existing records (previously rows 11-100) down by one row,
# which moves the former last record (row 100) to row 101 and grows the
# used range from A1:R100 to A1:R101, matching the target diff.
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value  = 2
$ws.Cells.Item(11, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(11, 3).Value  = "Coquimbo"
$ws.Cells.Item(11, 4).Value  = 44552
$ws.Cells.Item(11, 5).Value  = 4
$ws.Cells.Item(11, 6).Value  = 100112024
$ws.Cells.Item(11, 7).Value  = "Choclo"
$ws.Cells.Item(11, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(11, 9).Value  = "Primera"
$ws.Cells.Item(11, 10).Value = 660
$ws.Cells.Item(11, 11).Value = 10500
$ws.Cells.Item(11, 12).Value = 11000
$ws.Cells.Item(11, 13).Value = 10750
$ws.Cells.Item(11, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(11, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 16).Value = 154
$ws.Cells.Item(11, 17).Value = 70
$ws.Cells.Item(11, 18).Value = "Hortaliza"
